$d = $word.ActiveDocument

# wdReplaceAll = 2 / wdFindContinue = 1
$wdReplaceAll = 2
$wdFindContinue = 1

# ------------------------------------------------------------------
# 1. "You will edit an sbatch script" -> "You will edit a sbatch script"
# ------------------------------------------------------------------
$d.Content.Find.Execute("You will edit an sbatch script", $true, $false, $false, $false, $false, `
                         $true, $wdFindContinue, $false, "You will edit a sbatch script", $wdReplaceAll) | Out-Null

# ------------------------------------------------------------------
# 2. Drop "the outdirectory path," from the list of paths to edit
# ------------------------------------------------------------------
$d.Content.Find.Execute("file path, the outdirectory path, the path to the transcriptome directory", `
                         $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, `
                         "file path, the path to the transcriptome directory", $wdReplaceAll) | Out-Null

# ------------------------------------------------------------------
# 3. "(you do not need to do so)" -> "(you do not need to do this step)"
# ------------------------------------------------------------------
$d.Content.Find.Execute("(you do not need to do so)", $true, $false, $false, $false, $false, `
                         $true, $wdFindContinue, $false, "(you do not need to do this step)", $wdReplaceAll) | Out-Null

# ------------------------------------------------------------------
# 4. Reposition the curl-command screenshot (floating picture, editId 132B031C)
#    Document order of Shapes.Item(): writes land on the N-th shape in XML
#    document order; that anchor is the 7th anchored/floating shape.
# ------------------------------------------------------------------
$pic = $d.Shapes.Item(7)
$pic.Left = -847598 / 12700
$pic.Top = 7088251 / 12700

# ------------------------------------------------------------------
# 5. Remove the two blank "ListParagraph" paragraphs that sit between the
#    sbatch-script screenshot and "Check and see if the job is running".
#    Pattern: [image paragraph] [blank] [blank] [Check and see if ...]
# ------------------------------------------------------------------
for ($i = 1; $i -le ($d.Paragraphs.Count - 2); $i++) {
    $p0 = $d.Paragraphs.Item($i)
    $p1 = $d.Paragraphs.Item($i + 1)
    $p2 = $d.Paragraphs.Item($i + 2)
    if ($p0.Range.InlineShapes.Count -gt 0 -and `
        $p1.Range.Text.Trim().Length -eq 0 -and $p1.Range.InlineShapes.Count -eq 0 -and `
        $p2.Range.Text.Trim().Length -eq 0 -and $p2.Range.InlineShapes.Count -eq 0) {
        $p3 = $d.Paragraphs.Item($i + 3)
        if ($p3.Range.Text.StartsWith("Check and see if")) {
            $p1.Range.Delete()
            $p1again = $d.Paragraphs.Item($i + 1)
            $p1again.Range.Delete()
            break
        }
    }
}

Write-Output "done"
